$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# GPLIM-5135: Pooled Tube Upload template picked up a "Volume" value for
# the second (example) data row, which the template previously left blank.
$ws.Range("P2").Value = 61

# Reset the view: scroll back to the top-left of the sheet and leave the
# newly-filled cell selected (instead of the old P1:R2 block selection).
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("P2").Select()
